$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("N16").ClearContents()
$ws.Range("H21").Value = 58301.715
$ws.Range("I21").Value = 58301.715
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 58301.715
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = -57833.715
$ws.Range("N21").ClearContents()
$ws.Range("H23").Value = 58301.715
$ws.Range("I23").Value = 58301.715
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 58301.715
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = -58067.715
$ws.Range("N23").ClearContents()
$ws.Range("H32").Value = 866.2353000000001
$ws.Range("I32").Value = 366.6
$ws.Range("J32").Value = 1074.4166
$ws.Range("K32").Value = 366.6
$ws.Range("L32").Value = 1074.4166
$ws.Range("M32").Value = -40.60000000000002
$ws.Range("N32").Value = -1726.4166
$ws.Range("H98").Value = 1661.8948
$ws.Range("I98").Value = 1788.2222
$ws.Range("J98").Value = 1548.2
$ws.Range("K98").Value = 1788.2222
$ws.Range("L98").Value = 1548.2
$ws.Range("M98").Value = -290.2221999999999
$ws.Range("N98").Value = -4544.2
$ws.Range("H104").Value = 163.28572
$ws.Range("I104").Value = 157.25
$ws.Range("J104").Value = 171.33333
$ws.Range("K104").Value = 471.75
$ws.Range("L104").Value = 513.99999
$ws.Range("M104").Value = 1275.25
$ws.Range("N104").Value = -4007.99999
$ws.Range("H122").Value = 1661.8948
$ws.Range("I122").Value = 1788.2222
$ws.Range("J122").Value = 1548.2
$ws.Range("K122").Value = 5364.6666
$ws.Range("L122").Value = 4644.6
$ws.Range("M122").Value = -2914.6666
$ws.Range("N122").Value = -9544.6
$ws.Range("H127").Value = 2643.5144
$ws.Range("I127").Value = 500
$ws.Range("J127").Value = 2844.4688
$ws.Range("K127").Value = 1500
$ws.Range("L127").Value = 8533.4064
$ws.Range("M127").Value = 3460
$ws.Range("N127").Value = -18453.4064
$ws.Range("H129").Value = 29709608
$ws.Range("I129").Value = 478.14285
$ws.Range("J129").Value = 37411976
$ws.Range("K129").Value = 1434.42855
$ws.Range("L129").Value = 112235928
$ws.Range("M129").Value = 3565.57145
$ws.Range("N129").Value = -112245928
$ws.Range("H138").Value = 3118
$ws.Range("I138").Value = 1880.7693
$ws.Range("J138").Value = 3909.0164
$ws.Range("K138").Value = 5642.3079
$ws.Range("L138").Value = 11727.0492
$ws.Range("M138").Value = -502.3078999999998
$ws.Range("N138").Value = -22007.0492

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H12").Value = 25250
$ws.Range("I12").Value = 500
$ws.Range("K12").Value = 500
$ws.Range("M12").Value = -327
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").ClearContents()
$ws.Range("H32").Value = 5196.4194
$ws.Range("I32").Value = 4334.3335
$ws.Range("J32").Value = 17696.666
$ws.Range("K32").Value = 4334.3335
$ws.Range("L32").Value = 17696.666
$ws.Range("M32").Value = -4047.3335
$ws.Range("N32").Value = -18270.666
$ws.Range("H122").Value = 2511.0417
$ws.Range("I122").Value = 2327.1177
$ws.Range("K122").Value = 6981.353099999999
$ws.Range("M122").Value = -4531.353099999999
$ws.Range("H132").Value = 2308.8262
$ws.Range("I132").Value = 1652.4546
$ws.Range("J132").Value = 2910.5
$ws.Range("K132").Value = 4957.3638
$ws.Range("L132").Value = 8731.5
$ws.Range("M132").Value = -2427.3638
$ws.Range("N132").Value = -13791.5
$ws.Range("H141").Value = 51883.184
$ws.Range("J141").Value = 51883.184
$ws.Range("L141").Value = 51883.184
$ws.Range("N141").Value = -62243.184

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H51").Value = 33835.715
$ws.Range("J51").Value = 33835.715
$ws.Range("L51").Value = 33835.715
$ws.Range("N51").Value = -34817.715
$ws.Range("H107").Value = 1000
$ws.Range("I107").Value = 1000
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 1000
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 920
$ws.Range("N107").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 2246.2
$ws.Range("I22").Value = 2246.2
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 2246.2
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -1896.2
$ws.Range("N22").ClearContents()
$ws.Range("H58").Value = 2638.2727
$ws.Range("J58").Value = 2708.3333
$ws.Range("L58").Value = 2708.3333
$ws.Range("N58").Value = -3114.3333
$ws.Range("H136").Value = 2638.2727
$ws.Range("J136").Value = 2708.3333
$ws.Range("L136").Value = 8124.999899999999
$ws.Range("N136").Value = -13224.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("N31").ClearContents()
$ws.Range("H113").Value = 1814.5
$ws.Range("I113").Value = 1923.25
$ws.Range("J113").Value = 1705.75
$ws.Range("K113").Value = 5769.75
$ws.Range("L113").Value = 5117.25
$ws.Range("M113").Value = -3599.75
$ws.Range("N113").Value = -9457.25
$ws.Range("H121").Value = 35521.863
$ws.Range("I121").Value = 167066.5
$ws.Range("J121").Value = 1205.8695
$ws.Range("K121").Value = 501199.5
$ws.Range("L121").Value = 3617.6085
$ws.Range("M121").Value = -499889.5
$ws.Range("N121").Value = -6237.6085
$ws.Range("H131").Value = 833.2
$ws.Range("I131").Value = 362
$ws.Range("J131").Value = 858
$ws.Range("K131").Value = 1086
$ws.Range("L131").Value = 2574
$ws.Range("M131").Value = 3954
$ws.Range("N131").Value = -12654

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H27").Value = 5505.5
$ws.Range("J27").Value = 5505.5
$ws.Range("L27").Value = 5505.5
$ws.Range("N27").Value = -5837.5
$ws.Range("H132").Value = 2597.3635
$ws.Range("I132").Value = 2124.6667
$ws.Range("J132").Value = 4724.5
$ws.Range("K132").Value = 6374.000100000001
$ws.Range("L132").Value = 14173.5
$ws.Range("M132").Value = -3844.000100000001
$ws.Range("N132").Value = -19233.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("H22").Value = 1443.7693
$ws.Range("I22").Value = 957.8
$ws.Range("K22").Value = 957.8
$ws.Range("M22").Value = -662.8
$ws.Range("H27").Value = 1443.7693
$ws.Range("I27").Value = 957.8
$ws.Range("K27").Value = 957.8
$ws.Range("M27").Value = -850.8
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 0
$ws.Range("K28").Value = 0
$ws.Range("M28").ClearContents()
$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("M37").ClearContents()
$ws.Range("H45").Value = 17333.334
$ws.Range("I45").Value = 9500
$ws.Range("K45").Value = 9500
$ws.Range("M45").Value = -9093
$ws.Range("H53").Value = 7000
$ws.Range("I53").Value = 5000
$ws.Range("K53").Value = 5000
$ws.Range("M53").Value = -4482
$ws.Range("H132").Value = 11891.948
$ws.Range("I132").Value = 8864.529
$ws.Range("J132").Value = 32478.4
$ws.Range("K132").Value = 26593.587
$ws.Range("L132").Value = 97435.20000000001
$ws.Range("M132").Value = -24063.587
$ws.Range("N132").Value = -102495.2
$ws.Range("H133").Value = 36258.332
$ws.Range("J133").Value = 36258.332
$ws.Range("L133").Value = 36258.332
$ws.Range("N133").Value = -41318.332

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H22").Value = 21666.666
$ws.Range("J22").Value = 30000
$ws.Range("L22").Value = 30000
$ws.Range("N22").Value = -30586
